$wb = $excel.ActiveWorkbook

$data = @{
    "2025" = @{ B2 = 973.9537847600009; E2 = 28982.37596598056; I2 = 16175.28135478; L2 = 48524.529503538; M2 = 10590.587968015; N2 = 7155.076800480893; O2 = 6980.32583738884 }
    "2030" = @{ B2 = 5712.560177842886; E2 = 56106.05588781912; I2 = 44217.8984721661; L2 = 66966.57749858923; M2 = 21984.28023276101; N2 = 10592.74868329119; O2 = 12062.16737379073 }
    "2035" = @{ A2 = 2861.961401238371; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15130.68977466451; O2 = 14760.15955445609 }
    "2040" = @{ A2 = 2861.961401238371; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15235.09851993559; O2 = 14760.15955445609 }
    "2045" = @{ A2 = 6302.873118834019; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15770.60999844154; O2 = 17095.62619223827 }
    "2050" = @{ A2 = 6302.873118834019; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15770.60999844154; O2 = 17095.62619223827 }
}

foreach ($sheetName in $data.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $data[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}
